# unit 05 plant physiology fixes
$p = $ppt.ActivePresentation

# --- Slide 126 : "Q - 63" question text tweak ------------------------------
# Old:  ...வடிவம் என்ன?
# New:  ...வடிவம் எப்படி இருக்கும்?  (split into its own run, with the
#       separating space isolated into its own run too)
$s126 = $p.Slides.Item(126)
$sh126 = $s126.Shapes.Item(2)
$tr126 = $sh126.TextFrame.TextRange
$para126 = $tr126.Paragraphs(2)
$tail126 = $para126.Characters(45, 5)
$tail126.Text = "எப்படி இருக்கும்?"
$space126 = $para126.Characters(44, 1)
$space126.Text = " "
$sh126.Height = 285.9656

# --- Slide 127 : same question, Latha-styled version ------------------------
$s127 = $p.Slides.Item(127)
$sh127 = $s127.Shapes.Item(2)
$tr127 = $sh127.TextFrame.TextRange
$para127 = $tr127.Paragraphs(2)
$tail127 = $para127.Characters(45, 5)
$tail127.Text = "எப்படி இருக்கும்?"
$sh127.Height = 162.37031

# --- Slide 86 : picture repositioned / resized ------------------------------
$s86 = $p.Slides.Item(86)
$sh86 = $s86.Shapes.Item(4)
$sh86.Left = 233.3113
$sh86.Top = 213.5638
$sh86.Width = 354.4439
$sh86.Height = 195.41575

# --- Slide 87 : picture repositioned / resized ------------------------------
$s87 = $p.Slides.Item(87)
$sh87 = $s87.Shapes.Item(4)
$sh87.Left = 202.4543
$sh87.Top = 210.8572
$sh87.Width = 344.8928
$sh87.Height = 190.14993
